# "fix arima models and add sarima"
# Updates MAPE figures for the naive/ETS/arima family of models (rows 2-18),
# clears/re-applies the "best in row" yellow highlight (style index 2) where
# it moved, and appends two new model rows: arima1_0_2 and sarima212_001.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Recomputed MAPE values for existing rows (B:F = Winter/Spring/Summer/Autumn/year) ----
$ws.Range("B2").Value = 1.68812582128777
$ws.Range("C2").Value = 1.5454191152182
$ws.Range("D2").Value = 1.90362277850401
$ws.Range("E2").Value = 1.86316471821991
$ws.Range("F2").Value = 1.75002799016963

$ws.Range("B3").Value = 2.60986850849955
$ws.Range("C3").Value = 2.65830300752121
$ws.Range("D3").Value = 2.27943147165294
$ws.Range("E3").Value = 2.38277841681249
$ws.Range("F3").Value = 2.48234612320976

$ws.Range("B4").Value = 3.021852966882
$ws.Range("C4").Value = 5.02894907406369
$ws.Range("D4").Value = 3.72308847500766
$ws.Range("E4").Value = 4.54907834751226
$ws.Range("F4").Value = 4.08380650517374

$ws.Range("B5").Value = 1.90007191817192
$ws.Range("C5").Value = 1.85653146986307
$ws.Range("D5").Value = 1.93831198986972
$ws.Range("E5").Value = 1.89248542042924
$ws.Range("F5").Value = 1.89684891969641

$ws.Range("B6").Value = 1.68583697169122
$ws.Range("C6").Value = 1.63003563498664
$ws.Range("D6").Value = 1.83211784770186
$ws.Range("E6").Value = 1.77835072547725
$ws.Range("F6").Value = 1.73164509538667

$ws.Range("B7").Value = 1.68659209485817
$ws.Range("C7").Value = 1.63074154389817
$ws.Range("D7").Value = 1.83247531511126
$ws.Range("E7").Value = 1.77817542928494
$ws.Range("F7").Value = 1.73205608667158

$ws.Range("B8").Value = 1.60037788786339
$ws.Range("C8").Value = 1.50908721558408
$ws.Range("D8").Value = 1.53544442950787
$ws.Range("E8").Value = 1.54401427815876
$ws.Range("F8").Value = 1.5470216402532

$ws.Range("B9").Value = 2.9641727849579
$ws.Range("C9").Value = 2.91160280305163
$ws.Range("D9").Value = 3.56392728246669
$ws.Range("E9").Value = 3.39870215716373
$ws.Range("F9").Value = 3.21009111137695

$ws.Range("B10").Value = 2.9641727849579
$ws.Range("C10").Value = 2.91160280305163
$ws.Range("D10").Value = 3.56392728246669
$ws.Range("E10").Value = 3.39870215716373
# F10 stays an empty cell (no "year" MAPE figure for TBATS_rf)

$ws.Range("B11").Value = 1.58682330677456
$ws.Range("C11").Value = 1.51644749546971
$ws.Range("D11").Value = 1.53105280382889
$ws.Range("E11").Value = 1.56600924919838
$ws.Range("F11").Value = 1.54988886039842

$ws.Range("B12").Value = 4.92417787800037
$ws.Range("C12").Value = 5.45707803495353
$ws.Range("D12").Value = 7.19205840154032
$ws.Range("E12").Value = 5.29875511305917
$ws.Range("F12").Value = 5.72242234036895

$ws.Range("B13").Value = 5.37648976328945
$ws.Range("C13").Value = 5.52113108414505
$ws.Range("D13").Value = 8.87974220640168
$ws.Range("E13").Value = 7.08851758051864
$ws.Range("F13").Value = 6.72095148053404

$ws.Range("B14").Value = 1.53545441453833
$ws.Range("C14").Value = 1.67417829942883
$ws.Range("D14").Value = 1.81883647665146
$ws.Range("E14").Value = 1.89296311410148
$ws.Range("F14").Value = 1.73071307058369

$ws.Range("B15").Value = 1.3904171092013
$ws.Range("C15").Value = 1.47777523701746
$ws.Range("D15").Value = 1.56734222070151
$ws.Range("E15").Value = 1.51882112882175
$ws.Range("F15").Value = 1.48890910318917

$ws.Range("B16").Value = 1.39336301206938
$ws.Range("C16").Value = 1.4854102504412
$ws.Range("D16").Value = 1.54860620494541
$ws.Range("E16").Value = 1.51628653708042
$ws.Range("F16").Value = 1.48621324586092

$ws.Range("B17").Value = 2.77957147377093
$ws.Range("C17").Value = 12.4525138471674
$ws.Range("D17").Value = 15.542948420251
$ws.Range("E17").Value = 11.8730955924038
$ws.Range("F17").Value = 10.6910682514022

$ws.Range("B18").Value = 4.97407663646201
$ws.Range("C18").Value = 5.45597479967069
$ws.Range("D18").Value = 7.20118431619643
$ws.Range("E18").Value = 5.28056878297033
$ws.Range("F18").Value = 5.73226903833074

# ---- Re-apply the yellow "lowest MAPE in row" highlight where it moved ----
# Row 8 (holt-winters): highlight no longer applies to D:F
$ws.Range("D8:F8").Style = "Normal"
# Row 11 (ETS): highlight moves from C11 to D11
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Interior.Color = 65535
# Row 14 (arx): highlight no longer applies to B14
$ws.Range("B14").Style = "Normal"
# Row 15 (arima1_1_2): highlight now applies to B15:C15
$ws.Range("B15:C15").Interior.Color = 65535
# Row 16 (arima2_0_2): highlight now applies to E16:F16
$ws.Range("E16:F16").Interior.Color = 65535

# ---- Append two new model rows, copying the label formatting from A18 ----
$ws.Range("A18").Copy($ws.Range("A19:A20"))

$ws.Range("A19").Value = "arima1_0_2"
$ws.Range("B19").Value = 4.34502961991106
$ws.Range("C19").Value = 5.75598926276295
$ws.Range("D19").Value = 7.50972418272067
$ws.Range("E19").Value = 5.42101109565839
$ws.Range("F19").Value = 5.76465889819797

$ws.Range("A20").Value = "sarima212_001"
$ws.Range("B20").Value = 2.73284441964206
$ws.Range("C20").Value = 12.3331045517984
$ws.Range("D20").Value = 15.4148542117711
$ws.Range("E20").Value = 11.78061540627
$ws.Range("F20").Value = 5.76465889819797
